$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap LastName / FirstName labels (B1 <-> C1) ---
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"

# --- CNE numbers: 19000031..19000040 -> 18000001..18000010 ---
$ws.Range("A2").Value = 18000001
$ws.Range("A3").Value = 18000002
$ws.Range("A4").Value = 18000003
$ws.Range("A5").Value = 18000004
$ws.Range("A6").Value = 18000005
$ws.Range("A7").Value = 18000006
$ws.Range("A8").Value = 18000007
$ws.Range("A9").Value = 18000008
$ws.Range("A10").Value = 18000009
$ws.Range("A11").Value = 18000010

# --- Re-stamp formatting on the header row + the whole data block so it
#     picks up a freshly minted (default-looking) style record, matching
#     the workbook having been re-saved after an import/paste pass. ---
$ws.Range("A1:C1").Style = "Normal"
$ws.Range("A1:C11").HorizontalAlignment = 1

# --- Selection / view state left behind by the edit session ---
$ws.Range("D10").Select() | Out-Null
